$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.709.07'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '2.038.14'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.12'
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.604'
$ws.Range("E6").Value = '  -1.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.34'
$ws.Range("E7").Value = '  -0.63%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -2.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0836'
$ws.Range("E10").Value = '  +3.37%  '
$ws.Range("E11").Value = '  -0.20%  '
$ws.Range("D12").Value = '2.338.85'
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.42'
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.03'
$ws.Range("E14").Value = '  +0.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.45'
$ws.Range("E15").Value = '  +4.27%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.771'
$ws.Range("E16").Value = '  +2.73%  '
$ws.Range("D17").Value = '2.028.63'
$ws.Range("E17").Value = '  -0.09%  '
$ws.Range("D18").Value = '37.678.66'
$ws.Range("E18").Value = '  -0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.94'
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.31'
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("D21").Value = '0.0₃0823'
$ws.Range("E21").Value = '  +0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '223.81'
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -0.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.28'
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.00'
$ws.Range("E26").Value = '  +1.66%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.37'
$ws.Range("E27").Value = '  +2.01%  '
$ws.Range("E28").Value = '  -0.81%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.78'
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("E30").Value = '  -0.47%  '
$ws.Range("E31").Value = '  +0.30%  '
$ws.Range("E32").Value = '  +9.72%  '
$ws.Range("E33").Value = '  -1.21%  '
$ws.Range("E34").Value = '  +1.22%  '
$ws.Range("E35").Value = '  -0.43%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.54'
$ws.Range("E36").Value = '  +3.11%  '
$ws.Range("E37").Value = '  +4.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.42'
$ws.Range("E38").Value = '  +5.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.998'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.08'
$ws.Range("E40").Value = '  +8.79%  '
$ws.Range("D41").Value = '1.530.80'
$ws.Range("E41").Value = '  -0.48%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '97.13'
$ws.Range("E42").Value = '  +0.69%  '
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("E44").Value = '  +1.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.28'
$ws.Range("E45").Value = '  +9.95%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0908'
$ws.Range("E46").Value = '  -1.12%  '
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("E48").Value = '  +0.58%  '
$ws.Range("E49").Value = '  -0.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.05'
$ws.Range("E50").Value = '  -0.36%  '
$ws.Range("D51").Value = '2.228.73'
$ws.Range("E51").Value = '  +0.73%  '
